# Apply "Fruta / hortaliza, semanal" update: reorder Fecha/Volumen/Precio/Unidad/Origen
# data across rows 2-26 (columns D, M, N, O, P, Q, R, S, T) per the new weekly snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44238
$ws.Range("M2").Value = 60
$ws.Range("N2").Value = 15000
$ws.Range("O2").Value = 15000
$ws.Range("P2").Value = 15000
$ws.Range("Q2").Value = '$/caja 15 kilos granel'
$ws.Range("R2").Value = 'Provincia de Limarí'
$ws.Range("S2").Value = 1000
$ws.Range("T2").Value = 15

$ws.Range("D3").Value = 44616
$ws.Range("M3").Value = 70
$ws.Range("N3").Value = 14000
$ws.Range("O3").Value = 14000
$ws.Range("P3").Value = 14000
$ws.Range("Q3").Value = '$/caja 14 kilos empedrada'
$ws.Range("R3").Value = 'Provincia de Limarí'
$ws.Range("S3").Value = 1000
$ws.Range("T3").Value = 14

$ws.Range("D4").Value = 44260
$ws.Range("M4").Value = 56
$ws.Range("N4").Value = 13000
$ws.Range("O4").Value = 13000
$ws.Range("P4").Value = 13000
$ws.Range("Q4").Value = '$/caja 14 kilos empedrada'
$ws.Range("R4").Value = 'Provincia del Elquí'
$ws.Range("S4").Value = 929
$ws.Range("T4").Value = 14

$ws.Range("D5").Value = 44312
$ws.Range("M5").Value = 68
$ws.Range("N5").Value = 14000
$ws.Range("O5").Value = 14000
$ws.Range("P5").Value = 14000
$ws.Range("Q5").Value = '$/caja 14 kilos granel'
$ws.Range("R5").Value = 'Provincia de Limarí'
$ws.Range("S5").Value = 1000
$ws.Range("T5").Value = 14

$ws.Range("D6").Value = 44278
$ws.Range("M6").Value = 45
$ws.Range("N6").Value = 13000
$ws.Range("O6").Value = 13000
$ws.Range("P6").Value = 13000
$ws.Range("Q6").Value = '$/caja 14 kilos empedrada'
$ws.Range("R6").Value = 'Provincia del Elquí'
$ws.Range("S6").Value = 929
$ws.Range("T6").Value = 14

$ws.Range("D7").Value = 44242
$ws.Range("M7").Value = 45
$ws.Range("N7").Value = 12000
$ws.Range("O7").Value = 12000
$ws.Range("P7").Value = 12000
$ws.Range("Q7").Value = '$/caja 15 kilos granel'
$ws.Range("R7").Value = 'Provincia de Limarí'
$ws.Range("S7").Value = 800
$ws.Range("T7").Value = 15

$ws.Range("D8").Value = 44239
$ws.Range("M8").Value = 70
$ws.Range("N8").Value = 15000
$ws.Range("O8").Value = 15000
$ws.Range("P8").Value = 15000
$ws.Range("Q8").Value = '$/caja 15 kilos granel'
$ws.Range("R8").Value = 'Provincia de Limarí'
$ws.Range("S8").Value = 1000
$ws.Range("T8").Value = 15

$ws.Range("D9").Value = 44270
$ws.Range("M9").Value = 85
$ws.Range("N9").Value = 12000
$ws.Range("O9").Value = 12000
$ws.Range("P9").Value = 12000
$ws.Range("Q9").Value = '$/caja 14 kilos granel'
$ws.Range("R9").Value = 'Provincia del Elquí'
$ws.Range("S9").Value = 857
$ws.Range("T9").Value = 14

$ws.Range("D10").Value = 44592
$ws.Range("M10").Value = 54
$ws.Range("N10").Value = 20000
$ws.Range("O10").Value = 20000
$ws.Range("P10").Value = 20000
$ws.Range("Q10").Value = '$/caja 15 kilos empedrada'
$ws.Range("R10").Value = 'Provincia de Limarí'
$ws.Range("S10").Value = 1333
$ws.Range("T10").Value = 15

$ws.Range("D11").Value = 44588
$ws.Range("M11").Value = 85
$ws.Range("N11").Value = 19000
$ws.Range("O11").Value = 20000
$ws.Range("P11").Value = 19529
$ws.Range("Q11").Value = '$/caja 14 kilos granel'
$ws.Range("R11").Value = 'Provincia de Limarí'
$ws.Range("S11").Value = 1395
$ws.Range("T11").Value = 14

$ws.Range("D12").Value = 44252
$ws.Range("M12").Value = 60
$ws.Range("N12").Value = 14000
$ws.Range("O12").Value = 14000
$ws.Range("P12").Value = 14000
$ws.Range("Q12").Value = '$/caja 14 kilos empedrada'
$ws.Range("R12").Value = 'Provincia de Limarí'
$ws.Range("S12").Value = 1000
$ws.Range("T12").Value = 14

$ws.Range("D13").Value = 44322
$ws.Range("M13").Value = 50
$ws.Range("N13").Value = 14000
$ws.Range("O13").Value = 14000
$ws.Range("P13").Value = 14000
$ws.Range("Q13").Value = '$/caja 14 kilos granel'
$ws.Range("R13").Value = 'Provincia de Limarí'
$ws.Range("S13").Value = 1000
$ws.Range("T13").Value = 14

$ws.Range("D14").Value = 44316
$ws.Range("M14").Value = 48
$ws.Range("N14").Value = 14000
$ws.Range("O14").Value = 14000
$ws.Range("P14").Value = 14000
$ws.Range("Q14").Value = '$/caja 14 kilos granel'
$ws.Range("R14").Value = 'Provincia de Limarí'
$ws.Range("S14").Value = 1000
$ws.Range("T14").Value = 14

$ws.Range("D15").Value = 44614
$ws.Range("M15").Value = 54
$ws.Range("N15").Value = 14000
$ws.Range("O15").Value = 14000
$ws.Range("P15").Value = 14000
$ws.Range("Q15").Value = '$/caja 14 kilos granel'
$ws.Range("R15").Value = 'Provincia de Limarí'
$ws.Range("S15").Value = 1000
$ws.Range("T15").Value = 14

$ws.Range("D16").Value = 44585
$ws.Range("M16").Value = 50
$ws.Range("N16").Value = 22500
$ws.Range("O16").Value = 22500
$ws.Range("P16").Value = 22500
$ws.Range("Q16").Value = '$/caja 15 kilos empedrada'
$ws.Range("R16").Value = 'Provincia de Limarí'
$ws.Range("S16").Value = 1500
$ws.Range("T16").Value = 15

$ws.Range("D17").Value = 44271
$ws.Range("M17").Value = 50
$ws.Range("N17").Value = 12000
$ws.Range("O17").Value = 12000
$ws.Range("P17").Value = 12000
$ws.Range("Q17").Value = '$/caja 14 kilos granel'
$ws.Range("R17").Value = 'Provincia del Elquí'
$ws.Range("S17").Value = 857
$ws.Range("T17").Value = 14

$ws.Range("D18").Value = 44320
$ws.Range("M18").Value = 45
$ws.Range("N18").Value = 14000
$ws.Range("O18").Value = 14000
$ws.Range("P18").Value = 14000
$ws.Range("Q18").Value = '$/caja 14 kilos granel'
$ws.Range("R18").Value = 'Provincia de Limarí'
$ws.Range("S18").Value = 1000
$ws.Range("T18").Value = 14

$ws.Range("D19").Value = 44315
$ws.Range("M19").Value = 65
$ws.Range("N19").Value = 14000
$ws.Range("O19").Value = 14000
$ws.Range("P19").Value = 14000
$ws.Range("Q19").Value = '$/caja 14 kilos granel'
$ws.Range("R19").Value = 'Provincia de Limarí'
$ws.Range("S19").Value = 1000
$ws.Range("T19").Value = 14

$ws.Range("D20").Value = 44630
$ws.Range("M20").Value = 75
$ws.Range("N20").Value = 15000
$ws.Range("O20").Value = 15000
$ws.Range("P20").Value = 15000
$ws.Range("Q20").Value = '$/caja 14 kilos empedrada'
$ws.Range("R20").Value = 'Provincia de Limarí'
$ws.Range("S20").Value = 1071
$ws.Range("T20").Value = 14

$ws.Range("D21").Value = 44259
$ws.Range("M21").Value = 80
$ws.Range("N21").Value = 12000
$ws.Range("O21").Value = 12000
$ws.Range("P21").Value = 12000
$ws.Range("Q21").Value = '$/caja 15 kilos empedrada'
$ws.Range("R21").Value = 'Provincia de Limarí'
$ws.Range("S21").Value = 800
$ws.Range("T21").Value = 15

$ws.Range("D22").Value = 44313
$ws.Range("M22").Value = 36
$ws.Range("N22").Value = 14000
$ws.Range("O22").Value = 14000
$ws.Range("P22").Value = 14000
$ws.Range("Q22").Value = '$/caja 14 kilos granel'
$ws.Range("R22").Value = 'Provincia de Limarí'
$ws.Range("S22").Value = 1000
$ws.Range("T22").Value = 14

$ws.Range("D23").Value = 44314
$ws.Range("M23").Value = 56
$ws.Range("N23").Value = 14000
$ws.Range("O23").Value = 14000
$ws.Range("P23").Value = 14000
$ws.Range("Q23").Value = '$/caja 14 kilos granel'
$ws.Range("R23").Value = 'Provincia de Limarí'
$ws.Range("S23").Value = 1000
$ws.Range("T23").Value = 14

$ws.Range("D24").Value = 44627
$ws.Range("M24").Value = 56
$ws.Range("N24").Value = 17000
$ws.Range("O24").Value = 17000
$ws.Range("P24").Value = 17000
$ws.Range("Q24").Value = '$/caja 14 kilos empedrada'
$ws.Range("R24").Value = 'Provincia de Limarí'
$ws.Range("S24").Value = 1214
$ws.Range("T24").Value = 14

$ws.Range("D25").Value = 44245
$ws.Range("M25").Value = 50
$ws.Range("N25").Value = 15000
$ws.Range("O25").Value = 15000
$ws.Range("P25").Value = 15000
$ws.Range("Q25").Value = '$/caja 15 kilos granel'
$ws.Range("R25").Value = 'Provincia de Limarí'
$ws.Range("S25").Value = 1000
$ws.Range("T25").Value = 15

$ws.Range("D26").Value = 44323
$ws.Range("M26").Value = 60
$ws.Range("N26").Value = 14000
$ws.Range("O26").Value = 14000
$ws.Range("P26").Value = 14000
$ws.Range("Q26").Value = '$/caja 14 kilos granel'
$ws.Range("R26").Value = 'Provincia de Limarí'
$ws.Range("S26").Value = 1000
$ws.Range("T26").Value = 14
